$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2608695652173913
$ws.Range("C2").Value = 0.4347826086956522
$ws.Range("P2").Value = 0.2608695652173913
$ws.Range("S2").Value = 0.04347826086956522
$ws.Range("P3").Value = 0.9
$ws.Range("S3").Value = 0.1
$ws.Range("B6").Value = 0.1666666666666667
$ws.Range("J6").Value = 0.25
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.1666666666666667
$ws.Range("S6").Value = 0.25
$ws.Range("F7").Value = 0.1666666666666667
$ws.Range("J7").Value = 0.25
$ws.Range("O7").Value = 0.08333333333333333
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.25
$ws.Range("S7").Value = 0.08333333333333333
$ws.Range("B8").Value = 0.1111111111111111
$ws.Range("F8").Value = 0.1111111111111111
$ws.Range("J8").Value = 0.2592592592592592
$ws.Range("Q8").Value = 0.2592592592592592
$ws.Range("S8").Value = 0.2592592592592592
$ws.Range("J9").Value = 0.3
$ws.Range("Q9").Value = 0.3
$ws.Range("B10").Value = 0.09375
$ws.Range("D10").Value = 0.0234375
$ws.Range("F10").Value = 0.015625
$ws.Range("J10").Value = 0.0703125
$ws.Range("O10").Value = 0.015625
$ws.Range("Q10").Value = 0.2734375
$ws.Range("R10").Value = 0.1015625
$ws.Range("S10").Value = 0.40625
$ws.Range("G11").Value = 0.05555555555555555
$ws.Range("J11").Value = 0.1666666666666667
$ws.Range("K11").Value = 0.2222222222222222
$ws.Range("L11").Value = 0.5555555555555556
$ws.Range("G12").Value = 0.8
$ws.Range("J12").Value = 0.2
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.25
$ws.Range("G14").Value = 1
$ws.Range("I15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.4666666666666667
$ws.Range("F16").Value = 0.05555555555555555
$ws.Range("H16").Value = 0.1111111111111111
$ws.Range("I16").Value = 0.05555555555555555
$ws.Range("J16").Value = 0.3888888888888889
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.2222222222222222
$ws.Range("H17").Value = 0.125
$ws.Range("I17").Value = 0.0625
$ws.Range("J17").Value = 0.5416666666666666
$ws.Range("K17").Value = 0.0625
$ws.Range("M17").Value = 0.04166666666666666
$ws.Range("N17").Value = 0.02083333333333333
$ws.Range("O17").Value = 0.04166666666666666
$ws.Range("S17").Value = 0.1041666666666667
$ws.Range("H18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.6666666666666666
$ws.Range("O18").Value = 0.1111111111111111
$ws.Range("S18").Value = 0.1111111111111111
$ws.Range("H19").Value = 0.1683168316831683
$ws.Range("I19").Value = 0.04950495049504951
$ws.Range("J19").Value = 0.4752475247524752
$ws.Range("K19").Value = 0.07920792079207921
$ws.Range("M19").Value = 0.0198019801980198
$ws.Range("O19").Value = 0.06930693069306931
$ws.Range("S19").Value = 0.1386138613861386